$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 21.31228666666667
$ws.Cells.Item(2, 8).Value = 63.93686
$ws.Cells.Item(2, 9).Value = 0.9506775731819035
$ws.Cells.Item(2, 10).Value = 0.9506775731819034
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.04738633333333334
$ws.Cells.Item(2, 14).Value = 0.142159
$ws.Cells.Item(2, 15).Value = 0.05760194168856402
$ws.Cells.Item(2, 16).Value = 0.05760194168856402
$ws.Cells.Item(2, 17).Value = 1.009911120082222
$ws.Cells.Item(2, 18).Value = 9.089200080740001
$ws.Cells.Item(2, 19).Value = 0.05476087413504956
$ws.Cells.Item(2, 20).Value = 0.05476087413504956

$ws.Cells.Item(3, 7).Value = 21.31228666666667
$ws.Cells.Item(3, 8).Value = 63.93686
$ws.Cells.Item(3, 9).Value = 0.9506775731819035
$ws.Cells.Item(3, 10).Value = 0.9506775731819034
$ws.Cells.Item(3, 15).Value = 0.7659981644722047
$ws.Cells.Item(3, 16).Value = 0.7659981644722047
$ws.Cells.Item(3, 17).Value = 13.42993033890444
$ws.Cells.Item(3, 18).Value = 120.86937305014
$ws.Cells.Item(3, 19).Value = 0.7282172760622281
$ws.Cells.Item(3, 20).Value = 0.7282172760622281

$ws.Cells.Item(4, 7).Value = 21.31228666666667
$ws.Cells.Item(4, 8).Value = 63.93686
$ws.Cells.Item(4, 9).Value = 0.9506775731819035
$ws.Cells.Item(4, 10).Value = 0.9506775731819034
$ws.Cells.Item(4, 13).Value = 0.1451156666666667
$ws.Cells.Item(4, 14).Value = 0.435347
$ws.Cells.Item(4, 15).Value = 0.1763998938392313
$ws.Cells.Item(4, 16).Value = 0.1763998938392313
$ws.Cells.Item(4, 17).Value = 3.092746687824445
$ws.Cells.Item(4, 18).Value = 27.83472019042
$ws.Cells.Item(4, 19).Value = 0.1676994229846258
$ws.Cells.Item(4, 20).Value = 0.1676994229846258

$ws.Cells.Item(5, 7).Value = 0.7500946666666666
$ws.Cells.Item(5, 9).Value = 0.03345948693899053
$ws.Cells.Item(5, 10).Value = 0.03345948693899053
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.04738633333333334
$ws.Cells.Item(5, 14).Value = 0.142159
$ws.Cells.Item(5, 15).Value = 0.05760194168856402
$ws.Cells.Item(5, 16).Value = 0.05760194168856402
$ws.Cells.Item(5, 17).Value = 0.03554423590622222
$ws.Cells.Item(5, 18).Value = 0.319898123156
$ws.Cells.Item(5, 19).Value = 0.001927331415589002
$ws.Cells.Item(5, 20).Value = 0.001927331415589002

$ws.Cells.Item(6, 7).Value = 0.7500946666666666
$ws.Cells.Item(6, 9).Value = 0.03345948693899053
$ws.Cells.Item(6, 10).Value = 0.03345948693899053
$ws.Cells.Item(6, 15).Value = 0.7659981644722047
$ws.Cells.Item(6, 16).Value = 0.7659981644722047
$ws.Cells.Item(6, 17).Value = 0.4726719041684443
$ws.Cells.Item(6, 18).Value = 4.254047137515999
$ws.Cells.Item(6, 19).Value = 0.02562990557944845
$ws.Cells.Item(6, 20).Value = 0.02562990557944845

$ws.Cells.Item(7, 7).Value = 0.7500946666666666
$ws.Cells.Item(7, 9).Value = 0.03345948693899053
$ws.Cells.Item(7, 10).Value = 0.03345948693899053
$ws.Cells.Item(7, 13).Value = 0.1451156666666667
$ws.Cells.Item(7, 14).Value = 0.435347
$ws.Cells.Item(7, 15).Value = 0.1763998938392313
$ws.Cells.Item(7, 16).Value = 0.1763998938392313
$ws.Cells.Item(7, 18).Value = 0.979654388548
$ws.Cells.Item(7, 19).Value = 0.005902249943953075
$ws.Cells.Item(7, 20).Value = 0.005902249943953075

$ws.Cells.Item(8, 9).Value = 0.01586293987910606
$ws.Cells.Item(8, 10).Value = 0.01586293987910605
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.04738633333333334
$ws.Cells.Item(8, 14).Value = 0.142159
$ws.Cells.Item(8, 15).Value = 0.05760194168856402
$ws.Cells.Item(8, 16).Value = 0.05760194168856402
$ws.Cells.Item(8, 17).Value = 0.01685130672377778
$ws.Cells.Item(8, 18).Value = 0.151661760514
$ws.Cells.Item(8, 19).Value = 0.000913736137925464
$ws.Cells.Item(8, 20).Value = 0.0009137361379254638

$ws.Cells.Item(9, 9).Value = 0.01586293987910606
$ws.Cells.Item(9, 10).Value = 0.01586293987910605
$ws.Cells.Item(9, 15).Value = 0.7659981644722047
$ws.Cells.Item(9, 16).Value = 0.7659981644722047
$ws.Cells.Item(9, 19).Value = 0.01215098283052818
$ws.Cells.Item(9, 20).Value = 0.01215098283052818

$ws.Cells.Item(10, 9).Value = 0.01586293987910606
$ws.Cells.Item(10, 10).Value = 0.01586293987910605
$ws.Cells.Item(10, 13).Value = 0.1451156666666667
$ws.Cells.Item(10, 14).Value = 0.435347
$ws.Cells.Item(10, 15).Value = 0.1763998938392313
$ws.Cells.Item(10, 16).Value = 0.1763998938392313
$ws.Cells.Item(10, 17).Value = 0.05160535617355556
$ws.Cells.Item(10, 18).Value = 0.464448205562
$ws.Cells.Item(10, 19).Value = 0.002798220910652417
$ws.Cells.Item(10, 20).Value = 0.002798220910652416
